$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 4 swap order (23.10.0.8 driver now listed before 23.90.0.2),
# with updated Critical Minutes / Good Roaming Calculation values.
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 98.59999999999999

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.90.0.2"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = 372
$ws.Range("D4").Value = 98.7

# Totals row: Critical Minutes total updates to reflect new values.
$ws.Range("C5").Value = 380
